$wb = $excel.ActiveWorkbook

# Work on the "Geology" sheet (third sheet) which receives the new content.
$ws = $wb.Worksheets.Item("Geology")

# Add header/content cells.
$ws.Range("A1").Value = "https://www.google.com/"
$ws.Range("B1").Value = "instagram"
$ws.Range("C1").Value = "Instagram"
$ws.Range("B2").Value = "twitter"
$ws.Range("C2").Value = "Twitter. It's what's happening / Twitter"

# A1 becomes a hyperlink to google, styled like a normal hyperlink (blue).
$ws.Hyperlinks.Add($ws.Range("A1"), "https://www.google.com/", "", "", "https://www.google.com/")

# Column widths to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 32.22
$ws.Columns.Item(2).ColumnWidth = 31.68
$ws.Columns.Item(3).ColumnWidth = 35.15

# Make Geology the active sheet / tab and set its selection.
$ws.Activate()
$ws.Range("C20").Select()

# Political sheet selection moves too (no longer active tab).
$wsPolitical = $wb.Worksheets.Item("Political")
$wsPolitical.Range("G16").Select()

# "name" sheet's selection moves.
$wsName = $wb.Worksheets.Item("name")
$wsName.Range("D20").Select()

# Re-activate Geology last so it ends up the active tab.
$ws.Activate()
